$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.558.96"
$ws.Range("E2").Value = "  -1.12%  "

# Row 3
$ws.Range("D3").Value = "2.284.03"
$ws.Range("E3").Value = "  -0.87%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.10"
$ws.Range("E5").Value = "  +1.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.59"
$ws.Range("E6").Value = "  -2.82%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.504"
$ws.Range("E7").Value = "  -2.79%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  -3.10%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.77"
$ws.Range("E10").Value = "  -4.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").Value = "  -1.21%  "

# Row 12
$ws.Range("E12").Value = "  +1.56%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.95"
$ws.Range("E13").Value = "  +0.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.76"

# Row 15
$ws.Range("D15").Value = "2.639.36"
$ws.Range("E15").Value = "  -0.87%  "

# Row 16
$ws.Range("D16").Value = "2.264.76"
$ws.Range("E16").Value = "  -1.90%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.769"
$ws.Range("E17").Value = "  -1.73%  "

# Row 18
$ws.Range("D18").Value = "42.462.37"
$ws.Range("E18").Value = "  -1.24%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.86"
$ws.Range("E19").Value = "  +1.90%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  -2.22%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.96"
$ws.Range("E21").Value = "  -2.70%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.08"
$ws.Range("E22").Value = "  -1.88%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.95"
$ws.Range("E23").Value = "  -2.69%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  -0.35%  "

# Row 25
$ws.Range("E25").Value = "  +0.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  -1.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.67"
$ws.Range("E27").Value = "  -2.30%  "

# Row 28
$ws.Range("E28").Value = "  +16.88%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.49"
$ws.Range("E29").Value = "  -0.26%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.95"
$ws.Range("E30").Value = "  -1.59%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.45"
$ws.Range("E31").Value = "  -2.39%  "

# Row 32
$ws.Range("E32").Value = "  +0.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.79"
$ws.Range("E33").Value = "  +0.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.94"
$ws.Range("E34").Value = "  -1.78%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.43"
$ws.Range("E35").Value = "  -7.77%  "

# Row 36
$ws.Range("E36").Value = "  -2.51%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0683"
$ws.Range("E37").Value = "  -0.94%  "

# Row 38
$ws.Range("E38").Value = "  -0.76%  "

# Row 39
$ws.Range("E39").Value = "  -2.49%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.109"
$ws.Range("E40").Value = "  -2.02%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.66"
$ws.Range("E41").Value = "  -4.17%  "

# Row 42
$ws.Range("D42").Value = "1.989.32"
$ws.Range("E42").Value = "  -0.65%  "

# Row 43
$ws.Range("E43").Value = "  -3.83%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.35"
$ws.Range("E44").Value = "  +4.65%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.13"
$ws.Range("E45").Value = "  -0.83%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.01"
$ws.Range("E46").Value = "  -7.75%  "

# Row 47
$ws.Range("E47").Value = "  -2.25%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.93"
$ws.Range("E48").Value = "  +4.67%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.53"
$ws.Range("E49").Value = "  -0.38%  "

# Row 50
$ws.Range("D50").Value = "2.505.99"
$ws.Range("E50").Value = "  -0.89%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.13"
$ws.Range("E51").Value = "  +0.24%  "
